$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New boss rows appended after the existing data (rows 17-21)
$ws.Range("A17").Value = "WandBoss"
$ws.Range("B17").Value = 50
$ws.Range("C17").Value = 50
$ws.Range("D17").Value = 50

$ws.Range("A18").Value = "WegBoss"
$ws.Range("B18").Value = 150
$ws.Range("C18").Value = 250
$ws.Range("D18").Value = 350

$ws.Range("A19").Value = "SteinBoss"
$ws.Range("B19").Value = 75
$ws.Range("C19").Value = 75
$ws.Range("D19").Value = 75

$ws.Range("A20").Value = "RubinBoss"
$ws.Range("B20").Value = 75
$ws.Range("C20").Value = 150
$ws.Range("D20").Value = 75

$ws.Range("A21").Value = "BuschBoss"
$ws.Range("B21").Value = 75
$ws.Range("C21").Value = 75
$ws.Range("D21").Value = 150

# Update view: scroll back to top-left and move selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D22").Select() | Out-Null

$wb.Save() | Out-Null
